# Commit: "add the NA's under duplicate_image_filename"
#
# Column E ("duplicate_image_filename") is blank for every row in the
# stimuli table (rows 2-21). Fill those cells in with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
